{"js": "// Apply the dated worksheet refresh: update the header date and every\n// division problem's text, in document order. Each original string is\n// unique in the document at the time it is searched for, so a plain\n// search-and-replace per pair is unambiguous (one pair's replacement\n// text happens to equal an earlier pair's original text, but that\n// earlier occurrence has already been matched/replaced by then).\nconst replacements = [\n  [\"2024-05-27 Monday\", \"2024-05-28 Tuesday\"],\n  [\"411\u00f75=\", \"645\u00f79=\"],\n  [\"910\u00f73=\", \"809\u00f73=\"],\n  [\"328\u00f79=\", \"218\u00f73=\"],\n  [\"476\u00f78=\", \"170\u00f73=\"],\n  [\"898\u00f73=\", \"313\u00f73=\"],\n  [\"727\u00f76=\", \"378\u00f79=\"],\n  [\"587\u00f72=\", \"921\u00f77=\"],\n  [\"290\u00f78=\", \"362\u00f77=\"],\n  [\"607\u00f73=\", \"640\u00f76=\"],\n  [\"136\u00f79=\", \"929\u00f74=\"],\n  [\"218\u00f72=\", \"805\u00f72=\"],\n  [\"222\u00f79=\", \"632\u00f79=\"],\n  [\"397\u00f76=\", \"122\u00f78=\"],\n  [\"340\u00f77=\", \"556\u00f77=\"],\n  [\"808\u00f73=\", \"395\u00f73=\"],\n  [\"277\u00f75=\", \"789\u00f76=\"],\n  [\"194\u00f77=\", \"396\u00f73=\"],\n  [\"266\u00f75=\", \"132\u00f72=\"],\n  [\"868\u00f75=\", \"843\u00f78=\"],\n  [\"323\u00f75=\", \"584\u00f74=\"],\n  [\"998\u00f78=\", \"167\u00f72=\"],\n  [\"282\u00f79=\", \"633\u00f74=\"],\n  [\"719\u00f75=\", \"773\u00f74=\"],\n  [\"335\u00f75=\", \"537\u00f72=\"],\n  [\"164\u00f73=\", \"266\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the dated worksheet refresh: update the header date and every\n# division problem's text, in document order. Each pair's \"old\" text is\n# unique in the document at the moment it is searched for, so using\n# wdReplaceOne (replace exactly one match) per pair is unambiguous even\n# though one pair's \"new\" text happens to equal an earlier pair's \"old\"\n# text (that earlier occurrence has already been replaced by the time\n# the later pair runs).\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$pairs = @(\n  ,@(\"2024-05-27 Monday\", \"2024-05-28 Tuesday\")\n  ,@(\"411\u00f75=\", \"645\u00f79=\")\n  ,@(\"910\u00f73=\", \"809\u00f73=\")\n  ,@(\"328\u00f79=\", \"218\u00f73=\")\n  ,@(\"476\u00f78=\", \"170\u00f73=\")\n  ,@(\"898\u00f73=\", \"313\u00f73=\")\n  ,@(\"727\u00f76=\", \"378\u00f79=\")\n  ,@(\"587\u00f72=\", \"921\u00f77=\")\n  ,@(\"290\u00f78=\", \"362\u00f77=\")\n  ,@(\"607\u00f73=\", \"640\u00f76=\")\n  ,@(\"136\u00f79=\", \"929\u00f74=\")\n  ,@(\"218\u00f72=\", \"805\u00f72=\")\n  ,@(\"222\u00f79=\", \"632\u00f79=\")\n  ,@(\"397\u00f76=\", \"122\u00f78=\")\n  ,@(\"340\u00f77=\", \"556\u00f77=\")\n  ,@(\"808\u00f73=\", \"395\u00f73=\")\n  ,@(\"277\u00f75=\", \"789\u00f76=\")\n  ,@(\"194\u00f77=\", \"396\u00f73=\")\n  ,@(\"266\u00f75=\", \"132\u00f72=\")\n  ,@(\"868\u00f75=\", \"843\u00f78=\")\n  ,@(\"323\u00f75=\", \"584\u00f74=\")\n  ,@(\"998\u00f78=\", \"167\u00f72=\")\n  ,@(\"282\u00f79=\", \"633\u00f74=\")\n  ,@(\"719\u00f75=\", \"773\u00f74=\")\n  ,@(\"335\u00f75=\", \"537\u00f72=\")\n  ,@(\"164\u00f73=\", \"266\u00f75=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceOne)\n}\n"}
